$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("A2").Value = "Type"
$ws.Range("A1").Value = "Nb log"
$ws.Range("B1").Value = 0
$ws.Range("B2").Value = "message"

$ws.Range("C2").Select()
